$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update date / volume / price columns for existing data rows 2-28
$ws.Range("D2").Value = 44377
$ws.Range("J2").Value = 800
$ws.Range("K2").Value = 9000
$ws.Range("L2").Value = 10000
$ws.Range("M2").Value = 9500
$ws.Range("P2").Value = 380

$ws.Range("D3").Value = 44356
$ws.Range("J3").Value = 1000
$ws.Range("K3").Value = 11000
$ws.Range("L3").Value = 12000
$ws.Range("M3").Value = 11500
$ws.Range("P3").Value = 460

$ws.Range("D4").Value = 44349
$ws.Range("J4").Value = 600
$ws.Range("K4").Value = 10000
$ws.Range("L4").Value = 12000
$ws.Range("M4").Value = 11000
$ws.Range("P4").Value = 440

$ws.Range("D5").Value = 44364
$ws.Range("J5").Value = 700
$ws.Range("K5").Value = 11000
$ws.Range("L5").Value = 12000
$ws.Range("M5").Value = 11500
$ws.Range("P5").Value = 460

$ws.Range("D6").Value = 44336
$ws.Range("J6").Value = 1200
$ws.Range("K6").Value = 12000
$ws.Range("L6").Value = 13000
$ws.Range("M6").Value = 12500
$ws.Range("P6").Value = 500

$ws.Range("D7").Value = 44413
$ws.Range("J7").Value = 1200
$ws.Range("K7").Value = 10000
$ws.Range("L7").Value = 11000
$ws.Range("M7").Value = 10500
$ws.Range("P7").Value = 420

$ws.Range("D8").Value = 44308
$ws.Range("J8").Value = 400
$ws.Range("K8").Value = 11000
$ws.Range("L8").Value = 12000
$ws.Range("M8").Value = 11500
$ws.Range("P8").Value = 460

$ws.Range("D9").Value = 44363
$ws.Range("J9").Value = 900
$ws.Range("K9").Value = 11000
$ws.Range("L9").Value = 12000
$ws.Range("M9").Value = 11500
$ws.Range("P9").Value = 460

$ws.Range("D10").Value = 44328
$ws.Range("J10").Value = 900
$ws.Range("K10").Value = 11000
$ws.Range("L10").Value = 12000
$ws.Range("M10").Value = 11500
$ws.Range("P10").Value = 460

$ws.Range("D11").Value = 44343
$ws.Range("J11").Value = 500
$ws.Range("K11").Value = 9000
$ws.Range("L11").Value = 10000
$ws.Range("M11").Value = 9500
$ws.Range("P11").Value = 380

$ws.Range("D12").Value = 44406
$ws.Range("J12").Value = 800
$ws.Range("K12").Value = 10000
$ws.Range("L12").Value = 11000
$ws.Range("M12").Value = 10500
$ws.Range("P12").Value = 420

$ws.Range("D13").Value = 44385
$ws.Range("J13").Value = 600
$ws.Range("K13").Value = 8000
$ws.Range("L13").Value = 9000
$ws.Range("M13").Value = 8500
$ws.Range("P13").Value = 340

$ws.Range("D14").Value = 44371
$ws.Range("J14").Value = 500
$ws.Range("K14").Value = 10000
$ws.Range("L14").Value = 12000
$ws.Range("M14").Value = 11000
$ws.Range("P14").Value = 440

$ws.Range("D15").Value = 44419
$ws.Range("J15").Value = 1100
$ws.Range("K15").Value = 11000
$ws.Range("L15").Value = 12000
$ws.Range("M15").Value = 11500
$ws.Range("P15").Value = 460

$ws.Range("D16").Value = 44392
$ws.Range("J16").Value = 600
$ws.Range("K16").Value = 9000
$ws.Range("L16").Value = 10000
$ws.Range("M16").Value = 9500
$ws.Range("P16").Value = 380

$ws.Range("D17").Value = 44384
$ws.Range("J17").Value = 700
$ws.Range("K17").Value = 8000
$ws.Range("L17").Value = 9000
$ws.Range("M17").Value = 8500
$ws.Range("P17").Value = 340

$ws.Range("D18").Value = 44335
$ws.Range("J18").Value = 1000
$ws.Range("K18").Value = 12000
$ws.Range("L18").Value = 13000
$ws.Range("M18").Value = 12500
$ws.Range("P18").Value = 500

$ws.Range("D19").Value = 44426
$ws.Range("J19").Value = 500
$ws.Range("K19").Value = 11000
$ws.Range("L19").Value = 12000
$ws.Range("M19").Value = 11500
$ws.Range("P19").Value = 460

$ws.Range("D20").Value = 44434
$ws.Range("J20").Value = 600
$ws.Range("K20").Value = 10000
$ws.Range("L20").Value = 11000
$ws.Range("M20").Value = 10500
$ws.Range("P20").Value = 420

$ws.Range("D21").Value = 44398
$ws.Range("J21").Value = 400
$ws.Range("K21").Value = 9000
$ws.Range("L21").Value = 10000
$ws.Range("M21").Value = 9500
$ws.Range("P21").Value = 380

$ws.Range("D22").Value = 44420
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 10000
$ws.Range("L22").Value = 11000
$ws.Range("M22").Value = 10500
$ws.Range("P22").Value = 420

$ws.Range("D23").Value = 44427
$ws.Range("J23").Value = 360
$ws.Range("K23").Value = 10000
$ws.Range("L23").Value = 11000
$ws.Range("M23").Value = 10500
$ws.Range("P23").Value = 420

$ws.Range("D24").Value = 44441
$ws.Range("J24").Value = 1100
$ws.Range("K24").Value = 11000
$ws.Range("L24").Value = 12000
$ws.Range("M24").Value = 11500
$ws.Range("P24").Value = 460

$ws.Range("D25").Value = 44391
$ws.Range("J25").Value = 500
$ws.Range("K25").Value = 9000
$ws.Range("L25").Value = 10000
$ws.Range("M25").Value = 9500
$ws.Range("P25").Value = 380

$ws.Range("D26").Value = 44329
$ws.Range("J26").Value = 1000
$ws.Range("K26").Value = 12000
$ws.Range("L26").Value = 13000
$ws.Range("M26").Value = 12500
$ws.Range("P26").Value = 500

$ws.Range("D27").Value = 44435
$ws.Range("J27").Value = 600
$ws.Range("K27").Value = 10000
$ws.Range("L27").Value = 11000
$ws.Range("M27").Value = 10500
$ws.Range("P27").Value = 420

$ws.Range("D28").Value = 44412
$ws.Range("J28").Value = 1000
$ws.Range("K28").Value = 10000
$ws.Range("L28").Value = 11000
$ws.Range("M28").Value = 10500
$ws.Range("P28").Value = 420

# Append new data row 29 (same static Mercado/Categoria metadata as the other rows)
$ws.Range("A29").Value = 2
$ws.Range("B29").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C29").Value = "Coquimbo"
$ws.Range("D29").Value = 44399
$ws.Range("E29").Value = 4
$ws.Range("F29").Value = 100112026
$ws.Range("G29").Value = "Haba"
$ws.Range("H29").Value = "Sin especificar"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 500
$ws.Range("K29").Value = 9000
$ws.Range("L29").Value = 10000
$ws.Range("M29").Value = 9500
$ws.Range("N29").Value = "$/saco 25 kilos"
$ws.Range("O29").Value = "Provincia de Limarí"
$ws.Range("P29").Value = 380
$ws.Range("Q29").Value = 25
$ws.Range("R29").Value = "Hortaliza"

# Match the date number format used by the rest of column D
$ws.Range("D29").NumberFormat = "YYYY-MM-DD HH:MM:SS"
